$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.263.11"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "2.621.37"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "605.95"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "180.99"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "2.619.95"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +15.68%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "3.079.29"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").Value = "26.58"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("E16").Value = "  +7.64%  "
$ws.Range("D17").Value = "71.210.94"
$ws.Range("E17").Value = "  +4.56%  "
$ws.Range("D18").Value = "2.623.08"
$ws.Range("E18").Value = "  +4.81%  "
$ws.Range("D19").Value = "382.64"
$ws.Range("E19").Value = "  +8.82%  "
$ws.Range("E20").Value = "  +5.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.44%  "
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "72.31"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "4.45"
$ws.Range("E24").Value = "  +5.61%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +8.78%  "
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").Value = "2.755.90"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0955"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "544.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.34%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.32"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.83"
$ws.Range("E33").Value = "  +3.15%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "165.21"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "19.18"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  +7.14%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "18.98"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "1.38"
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.65"
$ws.Range("E41").Value = "  +9.11%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "5.05"
$ws.Range("E43").Value = "  +4.62%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.331"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "40.11"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "154.01"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.63"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.532"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "1.68"
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0264"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0753"
$ws.Range("E51").Value = "  +1.69%  "
